$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Sheet 1")

# New data rows feeding the SUM formula
$ws.Cells.Item(11, 3).Value = 1234.55
$ws.Cells.Item(12, 3).Value = 190

# Label + formula result row
$ws.Cells.Item(13, 2).Value = "Formula Result"
$ws.Range("C13").Formula = "=SUM(C10:C12)"

# Widen column B to fit the new label (stored width of 15 chars)
$ws.Columns.Item(2).ColumnWidth = 14.2

# Move the selection to the new formula cell, like the original edit did
$ws.Range("C13").Select() | Out-Null
